# Apply the cryptos.xlsx price/volume refresh described in the commit.
# Only the cells that actually change are touched; numeric-looking text
# values in column D are written with a leading apostrophe so Excel keeps
# them as text (matching the original inlineStr type) instead of silently
# coercing them to a Double and dropping formatting (e.g. "1.00" -> 1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.264.04"
$ws.Range("E2").Value = "  -1.69%  "

# Row 3
$ws.Range("D3").Value = "1.583.48"
$ws.Range("E3").Value = "  -1.08%  "

# Row 4
$ws.Range("E4").Value = "  -0.28%  "

# Row 5
$ws.Range("D5").Value = "'209.74"
$ws.Range("E5").Value = "  -0.66%  "

# Row 6
$ws.Range("D6").Value = "'0.506"
$ws.Range("E6").Value = "  -1.21%  "

# Row 7
$ws.Range("E7").Value = "  -0.26%  "

# Row 8
$ws.Range("E8").Value = "  -1.14%  "

# Row 9
$ws.Range("D9").Value = "'0.246"
$ws.Range("E9").Value = "  -0.20%  "

# Row 10
$ws.Range("E10").Value = "  -0.26%  "

# Row 11
$ws.Range("D11").Value = "'0.0847"
$ws.Range("E11").Value = "  +0.50%  "

# Row 12
$ws.Range("D12").Value = "1.804.63"
$ws.Range("E12").Value = "  -1.14%  "

# Row 13
$ws.Range("D13").Value = "1.598.68"
$ws.Range("E13").Value = "  -0.11%  "

# Row 14
$ws.Range("E14").Value = "  -0.36%  "

# Row 15
$ws.Range("E15").Value = "  -1.06%  "

# Row 16
$ws.Range("D16").Value = "'64.65"
$ws.Range("E16").Value = "  -0.77%  "

# Row 17
$ws.Range("D17").Value = "26.258.59"
$ws.Range("E17").Value = "  -1.60%  "

# Row 18
$ws.Range("D18").Value = "0.0₃0740"
$ws.Range("E18").Value = "  -0.18%  "

# Row 19
$ws.Range("D19").Value = "'7.24"
$ws.Range("E19").Value = "  +0.46%  "

# Row 20
$ws.Range("E20").Value = "  -0.19%  "

# Row 21
$ws.Range("D21").Value = "'206.78"
$ws.Range("E21").Value = "  -1.80%  "

# Row 22
$ws.Range("E22").Value = "  -0.85%  "

# Row 23
$ws.Range("D23").Value = "'2.21"
$ws.Range("E23").Value = "  -3.55%  "

# Row 24
$ws.Range("D24").Value = "'8.85"
$ws.Range("E24").Value = "  -1.24%  "

# Row 25
$ws.Range("D25").Value = "'144.54"
$ws.Range("E25").Value = "  +0.34%  "

# Row 26
$ws.Range("E26").Value = "  -0.28%  "

# Row 27
$ws.Range("E27").Value = "  -0.96%  "

# Row 28
$ws.Range("E28").Value = "  -0.87%  "

# Row 29
$ws.Range("D29").Value = "'15.27"
$ws.Range("E29").Value = "  -0.65%  "

# Row 30
$ws.Range("D30").Value = "'0.0504"
$ws.Range("E30").Value = "  -1.45%  "

# Row 31
$ws.Range("E31").Value = "  -1.17%  "

# Row 32
$ws.Range("E32").Value = "  -0.80%  "

# Row 34
$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Value = "1.285.38"
$ws.Range("E34").Value = "  -0.64%  "

# Row 35
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "'1.26"
$ws.Range("E35").Value = "  +6.64%  "

# Row 36
$ws.Range("E36").Value = "  -0.19%  "

# Row 37
$ws.Range("D37").Value = "'0.607"
$ws.Range("E37").Value = "  +0.24%  "

# Row 38
$ws.Range("E38").Value = "  -1.17%  "

# Row 39
$ws.Range("E39").Value = "  -1.48%  "

# Row 40
$ws.Range("D40").Value = "'0.818"
$ws.Range("E40").Value = "  -0.57%  "

# Row 41
$ws.Range("D41").Value = "'5.50"
$ws.Range("E41").Value = "  +1.62%  "

# Row 42
$ws.Range("D42").Value = "'0.769"
$ws.Range("E42").Value = "  -1.21%  "

# Row 43
$ws.Range("D43").Value = "'2.13"

# Row 44
$ws.Range("E44").Value = "  -1.43%  "

# Row 45
$ws.Range("D45").Value = "1.717.78"
$ws.Range("E45").Value = "  -1.23%  "

# Row 46
$ws.Range("D46").Value = "'88.94"
$ws.Range("E46").Value = "  -1.98%  "

# Row 47
$ws.Range("E47").Value = "  -0.38%  "

# Row 48
$ws.Range("E48").Value = "  +0.54%  "

# Row 49
$ws.Range("E49").Value = "  -1.56%  "

# Row 50
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  -0.21%  "

# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.39"
$ws.Range("E51").Value = "  -0.44%  "
